$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9176806807518005
$ws.Range("B1").Value = 3.041046619415283
$ws.Range("C1").Value = 4.276788234710693
$ws.Range("D1").Value = 3.008546829223633
$ws.Range("E1").Value = 1.383129954338074
